$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:D2").ClearContents()
$ws.Range("A2").Select()
